# Apply the commit's changes:
#  1. Refresh the cached "datetimeFigureOut" footer field text (27-Feb-19 -> 3/12/2019)
#     on the slide master, the notes master, and every slide layout.
#  2. Rename the Person-related labels to Place/Rating on slide 1
#     (UniquePersonList -> UniquePlaceList, Person -> Place, Phone -> Rating).

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sp = $shapes.Item($i)
        if ($sp.Name -like "Date Placeholder*") {
            $sp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "3/12/2019"

# Slide master footer date field.
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Notes master footer date field.
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# Every slide layout's footer date field.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Slide 1 text relabels (Person -> Place model rename).
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sp = $slide.Shapes.Item($i)
    if (-not $sp.HasTextFrame) { continue }
    if (-not $sp.TextFrame.HasText) { continue }
    $txt = $sp.TextFrame.TextRange.Text
    if ($txt -eq "UniquePersonList") {
        $sp.TextFrame.TextRange.Text = "UniquePlaceList"
    } elseif ($txt -eq "Person") {
        $sp.TextFrame.TextRange.Text = "Place"
    } elseif ($txt -eq "Phone") {
        $sp.TextFrame.TextRange.Text = "Rating"
    }
}
